$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.836.29'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.116.56'
$ws.Range('E3').Value = '  +6.63%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5323'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4379'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09016'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.08'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.47%  '
$ws.Range('E11').Value = '  +4.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.04'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.63%  '
$ws.Range('D13').Value = '2.116.91'
$ws.Range('E13').Value = '  +6.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.775'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.829'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001133'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06670'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.351'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.60%  '
$ws.Range('D23').Value = '30.902.65'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.02%  '
$ws.Range('D25').Value = '2.362.52'
$ws.Range('E25').Value = '  +6.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.274'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.587'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.185'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.55%  '
$ws.Range('E32').Value = '  +2.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.257'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.015'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.572'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +19.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02612'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '12.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.552'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.28%  '
$ws.Range('E39').Value = '  +4.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.502'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2289'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6875'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.77%  '
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6479'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.90%  '
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9993'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.233'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.278'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.42%  '
